$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row-group label cells (column A) ---
$ws.Range("A2").Value = "age_gr=30-39"
$ws.Range("A4").Value = "age_gr=40-48"
$ws.Range("A6").Value = "age_gr=49-57"
$ws.Range("A8").Value = "age_gr=>57"
$ws.Range("A10").Value = "educ_gr=low educ"
$ws.Range("A12").Value = "HHinc_gr=low inc"
$ws.Range("A14").Value = "expvol"
$ws.Range("A16").Value = "N"
$ws.Range("A17").Value = "R2"

# --- Update header row (unchanged, included for completeness) ---
$ws.Range("B1").Value = "incvar I"
$ws.Range("C1").Value = "incvar II"
$ws.Range("D1").Value = "incvar III"
$ws.Range("E1").Value = "inciqr I"
$ws.Range("F1").Value = "inciqr II"
$ws.Range("G1").Value = "inciqr III"
$ws.Range("C11").Value = "(0.02)"
$ws.Range("D11").Value = "(0.02)"
$ws.Range("F11").Value = "(0.03)"
$ws.Range("G11").Value = "(0.03)"

# --- Update text-valued data cells (non-numeric-looking strings) ---
$ws.Range("B2").Value = "-0.12***"
$ws.Range("C2").Value = "-0.12***"
$ws.Range("D2").Value = "-0.12***"
$ws.Range("E2").Value = "-0.46***"
$ws.Range("F2").Value = "-0.46***"
$ws.Range("G2").Value = "-0.44***"
$ws.Range("B3").Value = "(0.02)"
$ws.Range("C3").Value = "(0.02)"
$ws.Range("D3").Value = "(0.02)"
$ws.Range("E3").Value = "(0.05)"
$ws.Range("F3").Value = "(0.05)"
$ws.Range("G3").Value = "(0.05)"
$ws.Range("B4").Value = "-0.15***"
$ws.Range("C4").Value = "-0.15***"
$ws.Range("D4").Value = "-0.14***"
$ws.Range("E4").Value = "-0.91***"
$ws.Range("F4").Value = "-0.89***"
$ws.Range("G4").Value = "-0.88***"
$ws.Range("B5").Value = "(0.02)"
$ws.Range("C5").Value = "(0.02)"
$ws.Range("D5").Value = "(0.02)"
$ws.Range("E5").Value = "(0.05)"
$ws.Range("F5").Value = "(0.05)"
$ws.Range("G5").Value = "(0.05)"
$ws.Range("B6").Value = "-0.16***"
$ws.Range("C6").Value = "-0.17***"
$ws.Range("D6").Value = "-0.15***"
$ws.Range("E6").Value = "-1.06***"
$ws.Range("F6").Value = "-1.01***"
$ws.Range("G6").Value = "-0.99***"
$ws.Range("B7").Value = "(0.03)"
$ws.Range("C7").Value = "(0.03)"
$ws.Range("D7").Value = "(0.03)"
$ws.Range("E7").Value = "(0.05)"
$ws.Range("F7").Value = "(0.05)"
$ws.Range("G7").Value = "(0.05)"
$ws.Range("E8").Value = "-0.95***"
$ws.Range("F8").Value = "-0.90***"
$ws.Range("G8").Value = "-0.89***"
$ws.Range("B9").Value = "(0.04)"
$ws.Range("C9").Value = "(0.04)"
$ws.Range("D9").Value = "(0.04)"
$ws.Range("E9").Value = "(0.08)"
$ws.Range("F9").Value = "(0.08)"
$ws.Range("G9").Value = "(0.08)"
$ws.Range("C10").Value = "0.04**"
$ws.Range("F10").Value = "-0.25***"
$ws.Range("G10").Value = "-0.28***"
$ws.Range("D12").Value = "0.16***"
$ws.Range("G12").Value = "0.16***"
$ws.Range("D13").Value = "(0.02)"
$ws.Range("G13").Value = "(0.03)"
$ws.Range("B14").Value = "1.46***"
$ws.Range("C14").Value = "1.42***"
$ws.Range("D14").Value = "1.60***"
$ws.Range("E14").Value = "3.58***"
$ws.Range("F14").Value = "3.71***"
$ws.Range("G14").Value = "3.84***"
$ws.Range("B15").Value = "(0.49)"
$ws.Range("C15").Value = "(0.49)"
$ws.Range("D15").Value = "(0.49)"
$ws.Range("E15").Value = "(0.96)"
$ws.Range("F15").Value = "(0.95)"
$ws.Range("G15").Value = "(0.95)"

# --- Update numeric-looking data cells (force text storage) ---
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "-0.06"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "-0.07"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "-0.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "20602"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "20602"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20602"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "38815"
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "38815"
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "38815"
$ws.Range("G16").Style = "Normal"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "0.00"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "0.00"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.02"
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "0.03"
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "0.03"
$ws.Range("G17").Style = "Normal"

# --- Clear cells that no longer have values ---
$ws.Range("B10").ClearContents()
$ws.Range("B11").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("F13").ClearContents()
